$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "204204204"
$ws.Range("B7").Value = "Dan"
$ws.Range("C7").Value = "m"
$ws.Range("D7").Value = "1234"
$ws.Range("E7").Value = $false
